$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '67.993.99'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.86%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.517.06'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("E4").Value = '  +0.03%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '601.01'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.60%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '183.00'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +5.27%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +0.07%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.596'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +0.18%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.141'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +4.17%  '

$ws.Range("E10").Value = '  -2.32%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.435'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -0.69%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '4.126.93'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -0.12%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '32.25'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +11.69%  '

$ws.Range("E14").Value = '  -0.19%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.0000183'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +0.00%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '67.963.21'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +0.92%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '3.520.36'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +0.03%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '6.38'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +0.48%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '14.72'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +2.85%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '396.20'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.43%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '73.70'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.31%  '

$ws.Range("E23").Value = '  +0.83%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '5.71'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.24%  '

$ws.Range("E26").Value = '  +0.29%  '

$ws.Range("E27").Value = '  +0.94%  '

$ws.Range("E28").Value = '  -0.91%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.994'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -0.38%  '

$ws.Range("E30").Value = '  -0.08%  '

$ws.Range("E31").Value = '  -0.65%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '2.08'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -0.20%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '23.97'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -0.86%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '7.42'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +0.16%  '

$ws.Range("E35").Value = '  +0.08%  '

$ws.Range("E36").Value = '  +1.05%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '163.29'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.02%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '1.96'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +2.35%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.878'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -2.21%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '7.14'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +2.53%  '

$ws.Range("E41").Value = '  +0.55%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '27.73'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -0.06%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '2.68'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +2.02%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '26.78'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +1.14%  '

$ws.Range("E45").Value = '  -1.45%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.812.12'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +0.26%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '42.40'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -1.24%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.0304'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -1.21%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '343.36'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +0.62%  '

$ws.Range("E50").Value = '  -1.04%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '33.72'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.45%  '
